# Update NATMI TPM-derived ligand-receptor edge metrics for Efnb2-Ephb2
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 37.98277566666666
$ws.Range("H2").Value = 113.948327
$ws.Range("I2").Value = 0.697850645410475
$ws.Range("J2").Value = 0.6978506454104751
$ws.Range("M2").Value = 0.2087793333333333
$ws.Range("N2").Value = 0.626338
$ws.Range("O2").Value = 0.02275344108115409
$ws.Range("P2").Value = 0.02275344108115409
$ws.Range("Q2").Value = 7.930018581836221
$ws.Range("R2").Value = 71.37016723652599
$ws.Range("S2").Value = 0.0158785035437926
$ws.Range("T2").Value = 0.0158785035437926

$ws.Range("G3").Value = 37.98277566666666
$ws.Range("H3").Value = 113.948327
$ws.Range("I3").Value = 0.697850645410475
$ws.Range("J3").Value = 0.6978506454104751
$ws.Range("O3").Value = 0.9204452022087118
$ws.Range("P3").Value = 0.920445202208712
$ws.Range("Q3").Value = 320.7931288741518
$ws.Range("R3").Value = 2887.138159867366
$ws.Range("S3").Value = 0.6423332784263247
$ws.Range("T3").Value = 0.6423332784263249

$ws.Range("G4").Value = 37.98277566666666
$ws.Range("H4").Value = 113.948327
$ws.Range("I4").Value = 0.697850645410475
$ws.Range("J4").Value = 0.6978506454104751
$ws.Range("M4").Value = 0.05870933333333334
$ws.Range("N4").Value = 0.176128
$ws.Range("O4").Value = 0.006398331365399365
$ws.Range("P4").Value = 0.006398331365399365
$ws.Range("Q4").Value = 2.229943437539555
$ws.Range("R4").Value = 20.069490937856
$ws.Range("S4").Value = 0.004465079672894033
$ws.Range("T4").Value = 0.004465079672894034

$ws.Range("G5").Value = 37.98277566666666
$ws.Range("H5").Value = 113.948327
$ws.Range("I5").Value = 0.697850645410475
$ws.Range("J5").Value = 0.6978506454104751
$ws.Range("M5").Value = 0.4594193333333333
$ws.Range("N5").Value = 1.378258
$ws.Range("O5").Value = 0.05006899181852175
$ws.Range("P5").Value = 0.05006899181852175
$ws.Range("Q5").Value = 17.45002147492955
$ws.Range("R5").Value = 157.050193274366
$ws.Range("S5").Value = 0.0349406782556072
$ws.Range("T5").Value = 0.0349406782556072

$ws.Range("G6").Value = 37.98277566666666
$ws.Range("H6").Value = 113.948327
$ws.Range("I6").Value = 0.697850645410475
$ws.Range("J6").Value = 0.6978506454104751
$ws.Range("M6").Value = 0.003065
$ws.Range("N6").Value = 0.009195
$ws.Range("O6").Value = 0.0003340335262130221
$ws.Range("P6").Value = 0.0003340335262130222
$ws.Range("Q6").Value = 0.1164172074183333
$ws.Range("R6").Value = 1.047754866765
$ws.Range("S6").Value = 0.0002331055118564943
$ws.Range("T6").Value = 0.0002331055118564944

$ws.Range("I7").Value = 0.1779541659542351
$ws.Range("J7").Value = 0.1779541659542352
$ws.Range("M7").Value = 0.2087793333333333
$ws.Range("N7").Value = 0.626338
$ws.Range("O7").Value = 0.02275344108115409
$ws.Range("P7").Value = 0.02275344108115409
$ws.Range("Q7").Value = 2.022180321839778
$ws.Range("R7").Value = 18.199622896558
$ws.Range("S7").Value = 0.004049069630185606
$ws.Range("T7").Value = 0.004049069630185607

$ws.Range("I8").Value = 0.1779541659542351
$ws.Range("J8").Value = 0.1779541659542352
$ws.Range("O8").Value = 0.9204452022087118
$ws.Range("P8").Value = 0.920445202208712
$ws.Range("S8").Value = 0.1637970582656286
$ws.Range("T8").Value = 0.1637970582656287

$ws.Range("I9").Value = 0.1779541659542351
$ws.Range("J9").Value = 0.1779541659542352
$ws.Range("M9").Value = 0.05870933333333334
$ws.Range("N9").Value = 0.176128
$ws.Range("O9").Value = 0.006398331365399365
$ws.Range("P9").Value = 0.006398331365399365
$ws.Range("Q9").Value = 0.5686427707164445
$ws.Range("R9").Value = 5.117784936448
$ws.Range("S9").Value = 0.001138609721628466
$ws.Range("T9").Value = 0.001138609721628467

$ws.Range("I10").Value = 0.1779541659542351
$ws.Range("J10").Value = 0.1779541659542352
$ws.Range("M10").Value = 0.4594193333333333
$ws.Range("N10").Value = 1.378258
$ws.Range("O10").Value = 0.05006899181852175
$ws.Range("P10").Value = 0.05006899181852175
$ws.Range("Q10").Value = 4.449811772586445
$ws.Range("R10").Value = 40.048305953278
$ws.Range("S10").Value = 0.00890998567923446
$ws.Range("T10").Value = 0.008909985679234463

$ws.Range("I11").Value = 0.1779541659542351
$ws.Range("J11").Value = 0.1779541659542352
$ws.Range("M11").Value = 0.003065
$ws.Range("N11").Value = 0.009195
$ws.Range("O11").Value = 0.0003340335262130221
$ws.Range("P11").Value = 0.0003340335262130222
$ws.Range("Q11").Value = 0.02968676347166667
$ws.Range("R11").Value = 0.267180871245
$ws.Range("S11").Value = 0.00005944265755799049
$ws.Range("T11").Value = 0.00005944265755799051

$ws.Range("G12").Value = 0.5676613333333332
$ws.Range("H12").Value = 1.702984
$ws.Range("I12").Value = 0.01042953867610283
$ws.Range("J12").Value = 0.01042953867610283
$ws.Range("M12").Value = 0.2087793333333333
$ws.Range("N12").Value = 0.626338
$ws.Range("O12").Value = 0.02275344108115409
$ws.Range("P12").Value = 0.02275344108115409
$ws.Range("Q12").Value = 0.1185159547324444
$ws.Range("R12").Value = 1.066643592592
$ws.Range("S12").Value = 0.0002373078937703236
$ws.Range("T12").Value = 0.0002373078937703236

$ws.Range("G13").Value = 0.5676613333333332
$ws.Range("H13").Value = 1.702984
$ws.Range("I13").Value = 0.01042953867610283
$ws.Range("J13").Value = 0.01042953867610283
$ws.Range("O13").Value = 0.9204452022087118
$ws.Range("P13").Value = 0.920445202208712
$ws.Range("Q13").Value = 4.794327219763555
$ws.Range("R13").Value = 43.148944977872
$ws.Range("S13").Value = 0.00959981883566905
$ws.Range("T13").Value = 0.009599818835669053

$ws.Range("G14").Value = 0.5676613333333332
$ws.Range("H14").Value = 1.702984
$ws.Range("I14").Value = 0.01042953867610283
$ws.Range("J14").Value = 0.01042953867610283
$ws.Range("M14").Value = 0.05870933333333334
$ws.Range("N14").Value = 0.176128
$ws.Range("O14").Value = 0.006398331365399365
$ws.Range("P14").Value = 0.006398331365399365
$ws.Range("Q14").Value = 0.03332701843911111
$ws.Range("R14").Value = 0.299943165952
$ws.Range("S14").Value = 0.0000667316444379545
$ws.Range("T14").Value = 0.00006673164443795452

$ws.Range("G15").Value = 0.5676613333333332
$ws.Range("H15").Value = 1.702984
$ws.Range("I15").Value = 0.01042953867610283
$ws.Range("J15").Value = 0.01042953867610283
$ws.Range("M15").Value = 0.4594193333333333
$ws.Range("N15").Value = 1.378258
$ws.Range("O15").Value = 0.05006899181852175
$ws.Range("P15").Value = 0.05006899181852175
$ws.Range("Q15").Value = 0.2607945913191111
$ws.Range("R15").Value = 2.347151321872
$ws.Range("S15").Value = 0.0005221964866447487
$ws.Range("T15").Value = 0.0005221964866447489

$ws.Range("G16").Value = 0.5676613333333332
$ws.Range("H16").Value = 1.702984
$ws.Range("I16").Value = 0.01042953867610283
$ws.Range("J16").Value = 0.01042953867610283
$ws.Range("M16").Value = 0.003065
$ws.Range("N16").Value = 0.009195
$ws.Range("O16").Value = 0.0003340335262130221
$ws.Range("P16").Value = 0.0003340335262130222
$ws.Range("Q16").Value = 0.001739881986666666
$ws.Range("R16").Value = 0.01565893788
$ws.Range("S16").Value = 0.000003483815580753723
$ws.Range("T16").Value = 0.000003483815580753724

$ws.Range("G17").Value = 5.823095333333334
$ws.Range("H17").Value = 17.469286
$ws.Range("I17").Value = 0.1069866739681064
$ws.Range("J17").Value = 0.1069866739681064
$ws.Range("M17").Value = 0.2087793333333333
$ws.Range("N17").Value = 0.626338
$ws.Range("O17").Value = 0.02275344108115409
$ws.Range("P17").Value = 0.02275344108115409
$ws.Range("Q17").Value = 1.215741961629778
$ws.Range("R17").Value = 10.941677654668
$ws.Range("S17").Value = 0.002434314982601952
$ws.Range("T17").Value = 0.002434314982601952

$ws.Range("G18").Value = 5.823095333333334
$ws.Range("H18").Value = 17.469286
$ws.Range("I18").Value = 0.1069866739681064
$ws.Range("J18").Value = 0.1069866739681064
$ws.Range("O18").Value = 0.9204452022087118
$ws.Range("P18").Value = 0.920445202208712
$ws.Range("Q18").Value = 49.18042293975423
$ws.Range("R18").Value = 442.623806457788
$ws.Range("S18").Value = 0.09847537075421126
$ws.Range("T18").Value = 0.09847537075421127

$ws.Range("G19").Value = 5.823095333333334
$ws.Range("H19").Value = 17.469286
$ws.Range("I19").Value = 0.1069866739681064
$ws.Range("J19").Value = 0.1069866739681064
$ws.Range("M19").Value = 0.05870933333333334
$ws.Range("N19").Value = 0.176128
$ws.Range("O19").Value = 0.006398331365399365
$ws.Range("P19").Value = 0.006398331365399365
$ws.Range("Q19").Value = 0.3418700449564445
$ws.Range("R19").Value = 3.076830404608
$ws.Range("S19").Value = 0.0006845361917298912
$ws.Range("T19").Value = 0.0006845361917298912

$ws.Range("G20").Value = 5.823095333333334
$ws.Range("H20").Value = 17.469286
$ws.Range("I20").Value = 0.1069866739681064
$ws.Range("J20").Value = 0.1069866739681064
$ws.Range("M20").Value = 0.4594193333333333
$ws.Range("N20").Value = 1.378258
$ws.Range("O20").Value = 0.05006899181852175
$ws.Range("P20").Value = 0.05006899181852175
$ws.Range("Q20").Value = 2.675242575976445
$ws.Range("R20").Value = 24.077183183788
$ws.Range("S20").Value = 0.005356714903599975
$ws.Range("T20").Value = 0.005356714903599975

$ws.Range("G21").Value = 5.823095333333334
$ws.Range("H21").Value = 17.469286
$ws.Range("I21").Value = 0.1069866739681064
$ws.Range("J21").Value = 0.1069866739681064
$ws.Range("M21").Value = 0.003065
$ws.Range("N21").Value = 0.009195
$ws.Range("O21").Value = 0.0003340335262130221
$ws.Range("P21").Value = 0.0003340335262130222
$ws.Range("Q21").Value = 0.01784778719666667
$ws.Range("R21").Value = 0.16063008477
$ws.Range("S21").Value = 0.00003573713596336953
$ws.Range("T21").Value = 0.00003573713596336954

$ws.Range("G22").Value = 0.3689676666666666
$ws.Range("H22").Value = 1.106903
$ws.Range("I22").Value = 0.006778975991080511
$ws.Range("J22").Value = 0.006778975991080512
$ws.Range("M22").Value = 0.2087793333333333
$ws.Range("N22").Value = 0.626338
$ws.Range("O22").Value = 0.02275344108115409
$ws.Range("P22").Value = 0.02275344108115409
$ws.Range("Q22").Value = 0.07703282346822221
$ws.Range("R22").Value = 0.693295411214
$ws.Range("S22").Value = 0.0001542450308036086
$ws.Range("T22").Value = 0.0001542450308036086

$ws.Range("G23").Value = 0.3689676666666666
$ws.Range("H23").Value = 1.106903
$ws.Range("I23").Value = 0.006778975991080511
$ws.Range("J23").Value = 0.006778975991080512
$ws.Range("O23").Value = 0.9204452022087118
$ws.Range("P23").Value = 0.920445202208712
$ws.Range("Q23").Value = 3.116209654663778
$ws.Range("R23").Value = 28.045886891974
$ws.Range("S23").Value = 0.006239675926878103
$ws.Range("T23").Value = 0.006239675926878105

$ws.Range("G24").Value = 0.3689676666666666
$ws.Range("H24").Value = 1.106903
$ws.Range("I24").Value = 0.006778975991080511
$ws.Range("J24").Value = 0.006778975991080512
$ws.Range("M24").Value = 0.05870933333333334
$ws.Range("N24").Value = 0.176128
$ws.Range("O24").Value = 0.006398331365399365
$ws.Range("P24").Value = 0.006398331365399365
$ws.Range("Q24").Value = 0.02166184573155555
$ws.Range("R24").Value = 0.194956611584
$ws.Range("S24").Value = 0.00004337413470901968
$ws.Range("T24").Value = 0.00004337413470901968

$ws.Range("G25").Value = 0.3689676666666666
$ws.Range("H25").Value = 1.106903
$ws.Range("I25").Value = 0.006778975991080511
$ws.Range("J25").Value = 0.006778975991080512
$ws.Range("M25").Value = 0.4594193333333333
$ws.Range("N25").Value = 1.378258
$ws.Range("O25").Value = 0.05006899181852175
$ws.Range("P25").Value = 0.05006899181852175
$ws.Range("Q25").Value = 0.1695108794415555
$ws.Range("R25").Value = 1.525597914974
$ws.Range("S25").Value = 0.0003394164934353655
$ws.Range("T25").Value = 0.0003394164934353655

$ws.Range("G26").Value = 0.3689676666666666
$ws.Range("H26").Value = 1.106903
$ws.Range("I26").Value = 0.006778975991080511
$ws.Range("J26").Value = 0.006778975991080512
$ws.Range("M26").Value = 0.003065
$ws.Range("N26").Value = 0.009195
$ws.Range("O26").Value = 0.0003340335262130221
$ws.Range("P26").Value = 0.0003340335262130222
$ws.Range("Q26").Value = 0.001130885898333333
$ws.Range("R26").Value = 0.010177973085
$ws.Range("S26").Value = 0.00000226440525441404
$ws.Range("T26").Value = 0.00000226440525441404
